$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 678, shifting rows 678:737 down to 679:738
$ws.Rows("678:678").Insert()

# Populate the newly inserted row 678 with the new record's data
$ws.Range("A678").Value = 3
$ws.Range("B678").Value = "Femacal de La Calera"
$ws.Range("C678").Value = "Coquimbo"
$ws.Range("D678").Value = 45132
$ws.Range("E678").Value = 5
$ws.Range("F678").Value = 100112021
$ws.Range("G678").Value = "Ají"
$ws.Range("H678").Value = "Inferno"
$ws.Range("I678").Value = "Primera"
$ws.Range("J678").Value = 40
$ws.Range("K678").Value = 14000
$ws.Range("L678").Value = 14000
$ws.Range("M678").Value = 14000
$ws.Range("N678").Value = "$/caja 10 kilos"
$ws.Range("O678").Value = "Región de Arica y Parinacota"
$ws.Range("P678").Value = 1400
$ws.Range("Q678").Value = 10
$ws.Range("R678").Value = "Hortaliza"
